$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared string text updates (header volume/date labels) ---
$ws.Cells.Item(8, 1).Value = "Volume 29   Number  46"
$ws.Cells.Item(9, 3).Value = "Report Covering the Week  11/14/2022  Through  11/20/2022"

# --- Cell value / type updates ---
$ws.Cells.Item(14, 14).Value = -86.538461538461
$ws.Cells.Item(15, 3).NumberFormat = "@"
$ws.Cells.Item(15, 3).Value = "0"
$ws.Cells.Item(15, 3).NumberFormat = "General"
$ws.Cells.Item(15, 13).Value = -20
$ws.Cells.Item(15, 14).Value = -69.230769230769
$ws.Cells.Item(16, 3).Value = 5
$ws.Cells.Item(16, 4).Value = 6
$ws.Cells.Item(16, 5).Value = -16.666666666666
$ws.Cells.Item(16, 6).Value = 21
$ws.Cells.Item(16, 7).Value = 20
$ws.Cells.Item(16, 8).Value = 5
$ws.Cells.Item(16, 9).Value = 196
$ws.Cells.Item(16, 10).Value = 155
$ws.Cells.Item(16, 11).Value = 26.451612903225
$ws.Cells.Item(16, 12).Value = 29.801324503311
$ws.Cells.Item(16, 13).Value = -23.137254901960
$ws.Cells.Item(16, 14).Value = -75.950920245398
$ws.Cells.Item(17, 3).Value = 6
$ws.Cells.Item(17, 4).Value = 17
$ws.Cells.Item(17, 5).Value = -64.705882352941
$ws.Cells.Item(17, 6).Value = 41
$ws.Cells.Item(17, 7).Value = 40
$ws.Cells.Item(17, 8).Value = 2.5
$ws.Cells.Item(17, 9).Value = 462
$ws.Cells.Item(17, 10).Value = 468
$ws.Cells.Item(17, 11).Value = -1.282051282051
$ws.Cells.Item(17, 12).Value = 24.528301886792
$ws.Cells.Item(17, 13).Value = 62.105263157894
$ws.Cells.Item(17, 14).Value = -48.264277715565
$ws.Cells.Item(18, 3).Value = 1
$ws.Cells.Item(18, 5).Value = -75
$ws.Cells.Item(18, 6).Value = 13
$ws.Cells.Item(18, 7).Value = 19
$ws.Cells.Item(18, 8).Value = -31.578947368421
$ws.Cells.Item(18, 9).Value = 178
$ws.Cells.Item(18, 10).Value = 152
$ws.Cells.Item(18, 11).Value = 17.105263157894
$ws.Cells.Item(18, 12).Value = 60.360360360360
$ws.Cells.Item(18, 13).Value = 57.522123893805
$ws.Cells.Item(18, 14).Value = -79.610538373425
$ws.Cells.Item(19, 3).Value = 6
$ws.Cells.Item(19, 4).Value = 6
$ws.Cells.Item(19, 5).Value = 0
$ws.Cells.Item(19, 6).Value = 28
$ws.Cells.Item(19, 7).Value = 27
$ws.Cells.Item(19, 8).Value = 3.703703703703
$ws.Cells.Item(19, 9).Value = 311
$ws.Cells.Item(19, 10).Value = 283
$ws.Cells.Item(19, 11).Value = 9.893992932862
$ws.Cells.Item(19, 12).Value = 5.423728813559
$ws.Cells.Item(19, 13).Value = 14.338235294117
$ws.Cells.Item(19, 14).Value = -13.611111111111
$ws.Cells.Item(20, 4).Value = 1
$ws.Cells.Item(20, 5).Value = 0
$ws.Cells.Item(20, 6).Value = 8
$ws.Cells.Item(20, 7).Value = 9
$ws.Cells.Item(20, 8).Value = -11.111111111111
$ws.Cells.Item(20, 9).Value = 69
$ws.Cells.Item(20, 10).Value = 52
$ws.Cells.Item(20, 11).Value = 32.692307692307
$ws.Cells.Item(20, 12).Value = 23.214285714285
$ws.Cells.Item(20, 13).Value = 50
$ws.Cells.Item(20, 14).Value = -69.469026548672
$ws.Cells.Item(21, 3).Value = 19
$ws.Cells.Item(21, 4).Value = 34
$ws.Cells.Item(21, 5).Value = -44.117647058823
$ws.Cells.Item(21, 6).Value = 113
$ws.Cells.Item(21, 7).Value = 119
$ws.Cells.Item(21, 8).Value = -5.042016806722
$ws.Cells.Item(21, 9).Value = 1247
$ws.Cells.Item(21, 10).Value = 1144
$ws.Cells.Item(21, 11).Value = 9.003496503496
$ws.Cells.Item(21, 12).Value = 23.465346534653
$ws.Cells.Item(21, 13).Value = 23.710317460317
$ws.Cells.Item(21, 14).Value = -62.177737336973
$ws.Cells.Item(22, 3).Value = 2
$ws.Cells.Item(22, 3).NumberFormat = '#,##0'
$ws.Cells.Item(22, 6).Value = 4
$ws.Cells.Item(22, 9).Value = 19
$ws.Cells.Item(22, 11).Value = 58.333333333333
$ws.Cells.Item(22, 12).Value = 72.727272727272
$ws.Cells.Item(22, 13).Value = 111.111111111111
$ws.Cells.Item(23, 3).Value = 3
$ws.Cells.Item(23, 4).Value = 6
$ws.Cells.Item(23, 5).Value = -50
$ws.Cells.Item(23, 6).Value = 28
$ws.Cells.Item(23, 7).Value = 20
$ws.Cells.Item(23, 8).Value = 40
$ws.Cells.Item(23, 9).Value = 209
$ws.Cells.Item(23, 10).Value = 198
$ws.Cells.Item(23, 11).Value = 5.555555555555
$ws.Cells.Item(23, 12).Value = 17.415730337078
$ws.Cells.Item(23, 13).Value = 38.410596026490
$ws.Cells.Item(24, 3).Value = 13
$ws.Cells.Item(24, 4).Value = 16
$ws.Cells.Item(24, 5).Value = -18.75
$ws.Cells.Item(24, 6).Value = 47
$ws.Cells.Item(24, 7).Value = 50
$ws.Cells.Item(24, 8).Value = -6
$ws.Cells.Item(24, 9).Value = 732
$ws.Cells.Item(24, 10).Value = 681
$ws.Cells.Item(24, 11).Value = 7.488986784140
$ws.Cells.Item(24, 12).Value = 6.395348837209
$ws.Cells.Item(24, 13).Value = 27.526132404181
$ws.Cells.Item(25, 3).Value = 14
$ws.Cells.Item(25, 5).Value = 27.272727272727
$ws.Cells.Item(25, 6).Value = 58
$ws.Cells.Item(25, 7).Value = 38
$ws.Cells.Item(25, 8).Value = 52.631578947368
$ws.Cells.Item(25, 9).Value = 540
$ws.Cells.Item(25, 10).Value = 553
$ws.Cells.Item(25, 11).Value = -2.350813743218
$ws.Cells.Item(25, 12).Value = 15.138592750533
$ws.Cells.Item(25, 13).Value = -29.411764705882
$ws.Cells.Item(26, 3).NumberFormat = "@"
$ws.Cells.Item(26, 3).Value = "0"
$ws.Cells.Item(26, 3).NumberFormat = "General"
$ws.Cells.Item(27, 6).Value = 3
$ws.Cells.Item(27, 7).Value = 4
$ws.Cells.Item(27, 8).Value = -25
$ws.Cells.Item(27, 10).Value = 64
$ws.Cells.Item(27, 11).Value = 0
$ws.Cells.Item(27, 12).Value = 82.857142857142
$ws.Cells.Item(28, 4).Value = 4
$ws.Cells.Item(28, 4).NumberFormat = '#,##0'
$ws.Cells.Item(28, 5).Value = -100
$ws.Cells.Item(28, 5).NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Cells.Item(28, 7).Value = 6
$ws.Cells.Item(28, 8).Value = -66.666666666666
$ws.Cells.Item(28, 10).Value = 56
$ws.Cells.Item(28, 11).Value = -39.285714285714
$ws.Cells.Item(28, 14).Value = -71.428571428571
$ws.Cells.Item(29, 4).Value = 1
$ws.Cells.Item(29, 4).NumberFormat = '#,##0'
$ws.Cells.Item(29, 5).Value = -100
$ws.Cells.Item(29, 5).NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Cells.Item(29, 10).Value = 45
$ws.Cells.Item(29, 11).Value = -37.777777777777
$ws.Cells.Item(29, 14).Value = -74.074074074074
$ws.Cells.Item(30, 6).NumberFormat = "@"
$ws.Cells.Item(30, 6).Value = "0"
$ws.Cells.Item(30, 6).NumberFormat = "General"
